# Daily attendance processing - 2025-11-08 02:51:06
# Reorders the "Recorded By" (column G) names so that "System" is listed
# first whenever it appears alongside a recorded user's email address
# (dnasr281@gmail.com or backup@backdoor.com). Order of any remaining
# tokens (e.g. a trailing lowercase "system") is preserved.
#
# Note: this COM host shares scope between a function body and its caller,
# so loop-counter variable names must not collide across nested loops
# (e.g. an inner "$i" would clobber an outer "for ($i = ...)"). Each loop
# below therefore uses its own uniquely named counter.

function Test-ExactEquals($s1, $s2) {
    # Case-sensitive equality check (PowerShell's -eq/-ne/-cmatch operators
    # behave case-insensitively in this COM host), needed here because the
    # sheet mixes "System" and "system" tokens that must stay distinct.
    if ($s1 -eq $null -or $s2 -eq $null) { return $false }
    if ($s1.Length -ne $s2.Length) { return $false }
    for ($charIdx = 0; $charIdx -lt $s1.Length; $charIdx++) {
        $c1 = [int][char]$s1[$charIdx]
        $c2 = [int][char]$s2[$charIdx]
        if ($c1 -ne $c2) { return $false }
    }
    return $true
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 157; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = @($val -split ", ")

        if ($parts.Count -gt 1) {
            $sysIndex = -1
            for ($partIdx = 0; $partIdx -lt $parts.Count; $partIdx++) {
                if (Test-ExactEquals $parts[$partIdx].Trim() "System") {
                    $sysIndex = $partIdx
                    break
                }
            }

            if ($sysIndex -ge 0) {
                $first = $parts[0].Trim()
                $isKnownUser = (Test-ExactEquals $first "dnasr281@gmail.com") -or (Test-ExactEquals $first "backup@backdoor.com")

                if ($isKnownUser) {
                    $rest = New-Object System.Collections.ArrayList
                    for ($copyIdx = 0; $copyIdx -lt $parts.Count; $copyIdx++) {
                        if ($copyIdx -ne $sysIndex) {
                            [void]$rest.Add($parts[$copyIdx])
                        }
                    }
                    $newParts = @("System") + @($rest)
                    $newVal = [string]::Join(", ", $newParts)
                    $cell.Value = $newVal
                }
            }
        }
    }
}
